$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.633.44'
$ws.Range('E2').Value = '  +1.57%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.602.94'
$ws.Range('E3').Value = '  +1.52%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.74'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('E6').Value = '  +0.83%  '
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '27.47'
$ws.Range('E8').Value = '  +5.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.252'
$ws.Range('E9').Value = '  +1.33%  '
$ws.Range('E10').Value = '  +1.29%  '
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.831.04'
$ws.Range('E12').Value = '  +1.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.600.41'
$ws.Range('E13').Value = '  +1.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.543'
$ws.Range('E14').Value = '  +4.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.622.90'
$ws.Range('E15').Value = '  +1.47%  '
$ws.Range('E16').Value = '  +1.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.77'
$ws.Range('E17').Value = '  +2.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '241.28'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.67'
$ws.Range('E19').Value = '  +3.23%  '
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('E21').Value = '  +0.18%  '
$ws.Range('E22').Value = '  +0.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.32'
$ws.Range('E23').Value = '  +1.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.11'
$ws.Range('E24').Value = '  +1.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.98'
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.41'
$ws.Range('E26').Value = '  +1.85%  '
$ws.Range('E27').Value = '  +0.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.42'
$ws.Range('E28').Value = '  +1.19%  '
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0480'
$ws.Range('E30').Value = '  +2.53%  '
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('E32').Value = '  +0.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.19'
$ws.Range('E33').Value = '  +3.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.421.72'
$ws.Range('E34').Value = '  +0.19%  '
$ws.Range('E35').Value = '  +3.42%  '
$ws.Range('E36').Value = '  +4.60%  '
$ws.Range('E37').Value = '  -1.62%  '
$ws.Range('E38').Value = '  -0.17%  '
$ws.Range('E39').Value = '  +3.37%  '
$ws.Range('E40').Value = '  +3.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '55.99'
$ws.Range('E41').Value = '  +6.15%  '
$ws.Range('E42').Value = '  +0.65%  '
$ws.Range('E43').Value = '  +4.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.817'
$ws.Range('E44').Value = '  +3.58%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.986'
$ws.Range('E46').Value = '  +16.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '66.00'
$ws.Range('E47').Value = '  +2.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.33'
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.740.95'
$ws.Range('E49').Value = '  +1.40%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '86.37'
$ws.Range('E50').Value = '  +1.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0105'
$ws.Range('E51').Value = '  +1.90%  '

Write-Host "Updated cryptos list"
